$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Mohammed Shami, Kings XI Punjab): runs 2 -> 0, balls 2 -> 1
$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'1"

# Row 4 (Mohammed Shami, Kings XI Punjab): runs 0 -> 2, balls 1 -> 2
$ws.Range("C4").Value = "'2"
$ws.Range("D4").Value = "'2"
